$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Owner column (C) to "John" for every data row
$ws.Range("C2:C17").Value = "John"

# Set the Priority column (B) to "Critical" for every data row,
# except the rows whose priority should remain "Normal" or "High"
$ws.Range("B2:B4").Value = "Critical"
$ws.Range("B5").Value = "Normal"
$ws.Range("B6").Value = "Critical"
$ws.Range("B7").Value = "Critical"
$ws.Range("B8").Value = "Normal"
$ws.Range("B9:B15").Value = "Critical"
$ws.Range("B16").Value = "High"
$ws.Range("B17").Value = "Critical"
